$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 1.02; "C" = 1.044118180872593; "D" = 1.043145357024273; "E" = 1.047781559638126; "F" = 1.050528413494059; "I" = 1.038441739157099; "J" = 1.049184980738676; "K" = 1.0459200483339; "L" = 1.050543237551709; "M" = 1.053282439988449; "N" = 1.050674943754525 }
    3 = @{ "B" = 1.02; "C" = 1.046226395171011; "D" = 1.044740007509215; "E" = 1.049841109345477; "F" = 1.052644204640203; "I" = 1.039080072354664; "J" = 1.050934588122546; "K" = 1.047323625396632; "L" = 1.052411471262011; "M" = 1.055207341307905; "N" = 1.052427035781562 }
    4 = @{ "B" = 1.02; "C" = 1.047583123689852; "D" = 1.045764991823851; "E" = 1.051166763208405; "F" = 1.054006664732413; "I" = 1.039488231783777; "J" = 1.052059186050509; "K" = 1.048224510120399; "L" = 1.053612968247858; "M" = 1.056445929379291; "N" = 1.053553230767575 }
    5 = @{ "B" = 1.02; "C" = 1.048151752619412; "D" = 1.04619428180151; "E" = 1.051722425750708; "F" = 1.054577899773763; "I" = 1.039658666571405; "J" = 1.05253020088263; "K" = 1.048601516519301; "L" = 1.054116346346991; "M" = 1.056965001887149; "N" = 1.054024914494785 }
    6 = @{ "B" = 1.02; "C" = 1.048247127004992; "D" = 1.046266267706218; "E" = 1.051815628680786; "F" = 1.054673723209385; "I" = 1.039687215970475; "J" = 1.052609183729512; "K" = 1.048664717165461; "L" = 1.054200765216334; "M" = 1.05705206185702; "N" = 1.054104009506369 }
    7 = @{ "B" = 1.02; "C" = 1.047590728522164; "D" = 1.045770734318171; "E" = 1.05117419440297; "F" = 1.054014303617854; "I" = 1.039490513663058; "J" = 1.052065486669322; "K" = 1.048229554441429; "L" = 1.053619701167196; "M" = 1.056452871610024; "N" = 1.053559540333991 }
    8 = @{ "B" = 1.02; "C" = 1.044832222844307; "D" = 1.043685714961832; "E" = 1.048479068214722; "F" = 1.051244843344679; "I" = 1.038658486786227; "J" = 1.049777844793085; "K" = 1.046395926969798; "L" = 1.051176163818944; "M" = 1.053934429659831; "N" = 1.051268649743893 }
    9 = @{ "B" = 1.02; "C" = 1.039912759611939; "D" = 1.039957771015921; "E" = 1.04367454725389; "F" = 1.046312509637343; "I" = 1.037154300038514; "J" = 1.045687672574742; "K" = 1.043107491859093; "L" = 1.046812291251512; "M" = 1.049441809633731; "N" = 1.04717266901185 }
    10 = @{ "B" = 1.02; "C" = 1.036591240342679; "D" = 1.037434378865257; "E" = 1.040431966687876; "F" = 1.042986855234296; "I" = 1.036125024303177; "J" = 1.042919065936154; "K" = 1.040874860388977; "L" = 1.043861804744564; "M" = 1.046407674935259; "N" = 1.044400130634273 }
    11 = @{ "B" = 1.02; "C" = 1.035142497325875; "D" = 1.036332262681696; "E" = 1.039017979857524; "F" = 1.041537403082823; "I" = 1.035672850496819; "J" = 1.041709829458953; "K" = 1.039898134920864; "L" = 1.042573931198975; "M" = 1.045084099225069; "N" = 1.043189176902572 }
    12 = @{ "B" = 1.02; "C" = 1.034602743288098; "D" = 1.035921427873246; "E" = 1.038491224564911; "F" = 1.040997550469019; "I" = 1.035503900027239; "J" = 1.041259059261118; "K" = 1.039533800887574; "L" = 1.042093967679287; "M" = 1.04459095221554; "N" = 1.042737766559352 }
    13 = @{ "B" = 1.02; "C" = 1.03471859674529; "D" = 1.03600961999317; "E" = 1.038604285693088; "F" = 1.041113417573903; "I" = 1.035540185701462; "J" = 1.04135582434891; "K" = 1.039612021794556; "L" = 1.042196994138876; "M" = 1.044696803053206; "N" = 1.042834669064668 }
    14 = @{ "B" = 1.02; "C" = 1.035097914549701; "D" = 1.036298332943717; "E" = 1.038974469737545; "F" = 1.041492808782053; "I" = 1.035658905359446; "J" = 1.041672601658599; "K" = 1.039868050469332; "L" = 1.042534289953701; "M" = 1.045043366639303; "N" = 1.043151896234472 }
    15 = @{ "B" = 1.02; "C" = 1.035331407793386; "D" = 1.03647602396822; "E" = 1.039202347201084; "F" = 1.041726369164957; "I" = 1.035731920338691; "J" = 1.041867564550498; "K" = 1.040025593759344; "L" = 1.042741897039011; "M" = 1.045256694226292; "N" = 1.04334713599604 }
    16 = @{ "B" = 1.02; "C" = 1.03668716324626; "D" = 1.037507320015125; "E" = 1.040525595254721; "F" = 1.043082848185053; "I" = 1.036154895325005; "J" = 1.042999096078111; "K" = 1.040939469164661; "L" = 1.043947056108943; "M" = 1.046495306653939; "N" = 1.044480274428211 }
    17 = @{ "B" = 1.02; "C" = 1.037534746951478; "D" = 1.038151663367; "E" = 1.04135294355766; "F" = 1.043931176938158; "I" = 1.03641846615495; "J" = 1.043706059261809; "K" = 1.041510022206393; "L" = 1.044700234296201; "M" = 1.047269607982832; "N" = 1.045188241580719 }
    18 = @{ "B" = 1.02; "C" = 1.038028118090819; "D" = 1.038526586685194; "E" = 1.041834567068559; "F" = 1.044425086829593; "I" = 1.036571577533758; "J" = 1.044117417223339; "K" = 1.041841855499204; "L" = 1.045138560204192; "M" = 1.04772030464315; "N" = 1.045600183717738 }
    19 = @{ "B" = 1.02; "C" = 1.03819617501966; "D" = 1.038654272388935; "E" = 1.041998627676319; "F" = 1.044593345058711; "I" = 1.03662367907622; "J" = 1.044257511001683; "K" = 1.041954840187489; "L" = 1.045287851361558; "M" = 1.047873822383485; "N" = 1.045740476445315 }
    20 = @{ "B" = 1.02; "C" = 1.037443914067536; "D" = 1.038082625943931; "E" = 1.041264276000148; "F" = 1.043840253325071; "I" = 1.036390252276815; "J" = 1.043630312712999; "K" = 1.04144890688062; "L" = 1.044619528136277; "M" = 1.047186630339562; "N" = 1.045112387463121 }
    21 = @{ "B" = 1.02; "C" = 1.034986260199809; "D" = 1.03621335477315; "E" = 1.038865502583414; "F" = 1.041381128247341; "I" = 1.035623972935534; "J" = 1.041579363307003; "K" = 1.039792699056955; "L" = 1.042435008986732; "M" = 1.044941354394293; "N" = 1.043058525473722 }
    22 = @{ "B" = 1.02; "C" = 1.033431601286616; "D" = 1.035029605094553; "E" = 1.037348378816584; "F" = 1.039826499175029; "I" = 1.035136428488244; "J" = 1.040280537328251; "K" = 1.038742478210631; "L" = 1.041052294007216; "M" = 1.043520889843572; "N" = 1.041757855013103 }
    23 = @{ "B" = 1.02; "C" = 1.034256666068873; "D" = 1.035657948072766; "E" = 1.03815349593183; "F" = 1.040651456828984; "I" = 1.035395436510987; "J" = 1.040969966436188; "K" = 1.039300075422986; "L" = 1.04178618617516; "M" = 1.0442747510406; "N" = 1.042448263189452 }
    24 = @{ "B" = 1.02; "C" = 1.037484960656924; "D" = 1.038113823797872; "E" = 1.041304344000864; "F" = 1.043881340589297; "I" = 1.03640300284441; "J" = 1.043664542411295; "K" = 1.041476525232442; "L" = 1.044655998833693; "M" = 1.047224127268107; "N" = 1.045146665771515 }
    25 = @{ "B" = 1.02; "C" = 1.041191756574579; "D" = 1.040928108952249; "E" = 1.044923430502456; "F" = 1.047594059376801; "I" = 1.037547769186421; "J" = 1.046752294205976; "K" = 1.043964608781137; "L" = 1.0479475608015; "M" = 1.050609978817247; "N" = 1.048238802527915 }
}
foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
